# Update crypto price/volume data per the Dec 28 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.249.11'
$ws.Cells.Item(2, 5).Value = '  +0.38%  '

$ws.Cells.Item(3, 4).Value = '2.394.88'
$ws.Cells.Item(3, 5).Value = '  +5.53%  '

$ws.Cells.Item(4, 5).Value = '  -0.22%  '

$ws.Cells.Item(5, 4).Value = '''327.45'
$ws.Cells.Item(5, 5).Value = '  +7.92%  '

$ws.Cells.Item(6, 4).Value = '''105.94'
$ws.Cells.Item(6, 5).Value = '  -6.38%  '

$ws.Cells.Item(7, 4).Value = '''0.651'
$ws.Cells.Item(7, 5).Value = '  +2.64%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 4).Value = '''0.653'
$ws.Cells.Item(9, 5).Value = '  +5.96%  '

$ws.Cells.Item(10, 4).Value = '''42.22'
$ws.Cells.Item(10, 5).Value = '  -5.09%  '

$ws.Cells.Item(11, 5).Value = '  +1.22%  '

$ws.Cells.Item(12, 4).Value = '''8.74'
$ws.Cells.Item(12, 5).Value = '  -2.44%  '

$ws.Cells.Item(13, 4).Value = '''1.06'
$ws.Cells.Item(13, 5).Value = '  -0.29%  '

$ws.Cells.Item(14, 4).Value = '''17.07'
$ws.Cells.Item(14, 5).Value = '  +10.47%  '

$ws.Cells.Item(15, 5).Value = '  +1.95%  '

$ws.Cells.Item(16, 4).Value = '2.754.63'
$ws.Cells.Item(16, 5).Value = '  +5.60%  '

$ws.Cells.Item(17, 4).Value = '2.389.57'
$ws.Cells.Item(17, 5).Value = '  +6.45%  '

$ws.Cells.Item(18, 4).Value = '43.242.68'
$ws.Cells.Item(18, 5).Value = '  +0.67%  '

$ws.Cells.Item(19, 4).Value = '''7.75'
$ws.Cells.Item(19, 5).Value = '  +6.92%  '

$ws.Cells.Item(20, 5).Value = '  +1.12%  '

$ws.Cells.Item(21, 4).Value = '''76.98'
$ws.Cells.Item(21, 5).Value = '  +2.10%  '

$ws.Cells.Item(22, 4).Value = '''3.72'
$ws.Cells.Item(22, 5).Value = '  +3.96%  '

$ws.Cells.Item(23, 4).Value = '''274.52'
$ws.Cells.Item(23, 5).Value = '  +6.15%  '

$ws.Cells.Item(24, 5).Value = '  -0.36%  '

$ws.Cells.Item(25, 4).Value = '''9.66'
$ws.Cells.Item(25, 5).Value = '  +7.13%  '

$ws.Cells.Item(26, 5).Value = '  +1.43%  '

$ws.Cells.Item(27, 4).Value = '''0.999'
$ws.Cells.Item(27, 5).Value = '  -0.10%  '

$ws.Cells.Item(28, 4).Value = '''23.16'
$ws.Cells.Item(28, 5).Value = '  +3.40%  '

$ws.Cells.Item(29, 4).Value = '''176.35'
$ws.Cells.Item(29, 5).Value = '  +0.66%  '

$ws.Cells.Item(30, 5).Value = '  -1.92%  '

$ws.Cells.Item(31, 4).Value = '''37.38'
$ws.Cells.Item(31, 5).Value = '  -2.34%  '

$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32, 4).Value = '''0.0941'
$ws.Cells.Item(32, 5).Value = '  +5.04%  '

$ws.Cells.Item(33, 2).Value = 'WEMIXToken'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(33, 4).Value = '''3.16'
$ws.Cells.Item(33, 5).Value = '  -0.79%  '

$ws.Cells.Item(34, 5).Value = '  +4.09%  '

$ws.Cells.Item(35, 5).Value = '  +5.00%  '

$ws.Cells.Item(36, 4).Value = '''4.90'
$ws.Cells.Item(36, 5).Value = '  -3.72%  '

$ws.Cells.Item(37, 4).Value = '''4.15'
$ws.Cells.Item(37, 5).Value = '  -2.95%  '

$ws.Cells.Item(38, 5).Value = '  -3.38%  '

$ws.Cells.Item(39, 5).Value = '  +4.10%  '

$ws.Cells.Item(40, 4).Value = '''2.83'
$ws.Cells.Item(40, 5).Value = '  +15.31%  '

$ws.Cells.Item(41, 4).Value = '''1.59'
$ws.Cells.Item(41, 5).Value = '  +17.74%  '

$ws.Cells.Item(42, 4).Value = '''0.236'
$ws.Cells.Item(42, 5).Value = '  +1.08%  '

$ws.Cells.Item(43, 4).Value = '''70.32'
$ws.Cells.Item(43, 5).Value = '  -2.77%  '

$ws.Cells.Item(44, 4).Value = '''122.68'
$ws.Cells.Item(44, 5).Value = '  +13.69%  '

$ws.Cells.Item(45, 5).Value = '  +0.18%  '

$ws.Cells.Item(46, 2).Value = 'BitcoinSV'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(46, 4).Value = '''93.32'
$ws.Cells.Item(46, 5).Value = '  +43.68%  '

$ws.Cells.Item(47, 2).Value = 'Celestia'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(47, 4).Value = '''12.37'
$ws.Cells.Item(47, 5).Value = '  -1.94%  '

$ws.Cells.Item(48, 4).Value = '''9.35'
$ws.Cells.Item(48, 5).Value = '  +7.39%  '

$ws.Cells.Item(49, 5).Value = '  +0.15%  '

$ws.Cells.Item(50, 5).Value = '  +0.76%  '

$ws.Cells.Item(51, 4).Value = '1.593.59'
$ws.Cells.Item(51, 5).Value = '  +8.44%  '
